$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = "Data analysis skills"
$ws.Range("C6").Value = "Thoroughness"
$ws.Range("C11").Value = "Strategic planning"
$ws.Range("C16").Value = "Risk assessment"
